# Swap the theme applied to the Slide Master with the theme applied to the
# Notes Master. Before this edit:
#   ppt/theme/theme1.xml (Slide Master's theme) = "Integral" / "Red Violet"
#   ppt/theme/theme2.xml (Notes Master's theme) = "Office Theme" / "Office"
# After this edit the two colour schemes are exchanged, i.e. the deck's
# Slide Master now uses the stock "Office Theme" colours and the Notes
# Master uses the colours that used to belong to the Slide Master.
#
# PowerPoint's ColorScheme.Colors(n).RGB setter stores values using the
# classic Windows COLORREF byte order (0x00BBGGRR), so a little helper is
# used below to turn an "RRGGBB" hex string into the integer that has to be
# assigned so that the OOXML <a:srgbClr val="RRGGBB"/> ends up correct.
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

function SetColorScheme($colorScheme, $hexColors) {
    # Order matches the 12 ColorScheme.Colors() indices / clrScheme children:
    # 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
    for ($i = 0; $i -lt $hexColors.Count; $i++) {
        $colorScheme.Colors($i + 1).RGB = HexToComRgb($hexColors[$i])
    }
}

# Colours that used to belong to the Slide Master's theme ("Integral" / "Red Violet").
$integralRedViolet = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "454551", # dk2
    "D8D9DC", # lt2
    "E32D91", # accent1
    "C830CC", # accent2
    "4EA6DC", # accent3
    "4775E7", # accent4
    "8971E1", # accent5
    "D54773", # accent6
    "6B9F25", # hlink
    "8C8C8C"  # folHlink
)

# Colours that used to belong to the Notes Master's theme ("Office Theme" / "Office").
$officeTheme = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation

# Slide Master: Integral/Red Violet -> Office Theme
SetColorScheme $p.SlideMaster.ColorScheme $officeTheme

# Notes Master: Office Theme -> Integral/Red Violet
SetColorScheme $p.NotesMaster.ColorScheme $integralRedViolet
